# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") for each trade row with the recalculated value,
# and refreshes the dependent IP/I0/IF values (H2/I2/J2) for the first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values, keyed by row number.
$kValues = @{
    2  = 0
    3  = 2
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 3
    12 = 2
    13 = 1
    14 = 3
    15 = 0
    16 = 2
    17 = 1
    18 = 3
    19 = 0
    20 = 1
    21 = 0
    22 = 4
    23 = 1
    24 = 1
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 2
    36 = 2
    37 = 2
    38 = 1
    39 = 3
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    46 = 2
    47 = 2
    48 = 1
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 1
    59 = 2
    60 = 0
    61 = 0
    62 = 1
    63 = 1
    64 = 0
    65 = 0
    66 = 1
    67 = 1
    68 = 3
    69 = 3
    70 = 1
    71 = 1
    72 = 2
    73 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

# Row 2 also had its recalculated IP / I0 / IF (std/mean derived) values change.
$ws.Cells.Item(2, 8).Value = 2
$ws.Cells.Item(2, 9).Value = 8
$ws.Cells.Item(2, 10).Value = 9
